$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in column F, row 1, matching the style of the
# other header cells (bold/centered/bordered style already used by B1:E1).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Timestamps (as plain text, matching the source data) for rows 2 through 84.
$timeTaken = @(
    "2021-10-05 13:39:53.979386",
    "2021-10-05 13:39:53.979397",
    "2021-10-05 13:39:53.979400",
    "2021-10-05 13:39:53.979402",
    "2021-10-05 13:39:53.979405",
    "2021-10-05 13:39:53.979408",
    "2021-10-05 13:39:53.979410",
    "2021-10-05 13:39:53.979413",
    "2021-10-05 13:39:53.979416",
    "2021-10-05 13:39:53.979418",
    "2021-10-05 13:39:53.979421",
    "2021-10-05 13:39:53.979423",
    "2021-10-05 13:39:53.979426",
    "2021-10-05 13:39:53.979428",
    "2021-10-05 13:39:53.979430",
    "2021-10-05 13:39:53.979433",
    "2021-10-05 13:39:53.979436",
    "2021-10-05 13:39:53.979438",
    "2021-10-05 13:39:53.979441",
    "2021-10-05 13:39:53.979443",
    "2021-10-05 13:39:53.979446",
    "2021-10-05 13:39:53.979448",
    "2021-10-05 13:39:53.979451",
    "2021-10-05 13:39:53.979453",
    "2021-10-05 13:39:53.979456",
    "2021-10-05 13:39:53.979458",
    "2021-10-05 13:39:53.979461",
    "2021-10-05 13:39:53.979463",
    "2021-10-05 13:39:53.979466",
    "2021-10-05 13:39:53.979468",
    "2021-10-05 13:39:53.979471",
    "2021-10-05 13:39:53.979473",
    "2021-10-05 13:39:53.979476",
    "2021-10-05 13:39:53.979478",
    "2021-10-05 13:39:53.979481",
    "2021-10-05 13:39:53.979483",
    "2021-10-05 13:39:53.979486",
    "2021-10-05 13:39:53.979488",
    "2021-10-05 13:39:53.979491",
    "2021-10-05 13:39:53.979493",
    "2021-10-05 13:39:53.979496",
    "2021-10-05 13:39:53.979499",
    "2021-10-05 13:39:53.979501",
    "2021-10-05 13:39:53.979504",
    "2021-10-05 13:39:53.979506",
    "2021-10-05 13:39:53.979509",
    "2021-10-05 13:39:53.979511",
    "2021-10-05 13:39:53.979514",
    "2021-10-05 13:39:53.979516",
    "2021-10-05 13:39:53.979519",
    "2021-10-05 13:39:53.979521",
    "2021-10-05 13:39:53.979524",
    "2021-10-05 13:39:53.979527",
    "2021-10-05 13:39:53.979529",
    "2021-10-05 13:39:53.979532",
    "2021-10-05 13:39:53.979534",
    "2021-10-05 13:39:53.979537",
    "2021-10-05 13:39:53.979539",
    "2021-10-05 13:39:53.979541",
    "2021-10-05 13:39:53.979544",
    "2021-10-05 13:39:53.979546",
    "2021-10-05 13:39:53.979549",
    "2021-10-05 13:39:53.979551",
    "2021-10-05 13:39:53.979554",
    "2021-10-05 13:39:53.979557",
    "2021-10-05 13:39:53.979560",
    "2021-10-05 13:39:53.979563",
    "2021-10-05 13:39:53.979565",
    "2021-10-05 13:39:53.979568",
    "2021-10-05 13:39:53.979570",
    "2021-10-05 13:39:53.979573",
    "2021-10-05 13:39:53.979575",
    "2021-10-05 13:39:53.979578",
    "2021-10-05 13:39:53.979580",
    "2021-10-05 13:39:53.979583",
    "2021-10-05 13:39:53.979586",
    "2021-10-05 13:39:53.979590",
    "2021-10-05 13:39:53.979593",
    "2021-10-05 13:39:53.979596",
    "2021-10-05 13:39:53.979598",
    "2021-10-05 13:39:53.979601",
    "2021-10-05 13:39:53.979607",
    "2021-10-05 13:39:53.979609"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

# Update the worksheet dimension to reflect the new column F.
$ws.UsedRange | Out-Null
